# Commit: "Fixed POI packaging and upgraded to POI 3.15."
#
# This commit regenerated the expected-generation test fixture with a
# newer Apache POI release. Diffing the two OOXML payloads shows that
# every changed line is the *same* element with the *same* attribute
# name/value pairs, just re-emitted with the attributes (and the
# xmlns:* declarations on the <w:document> root) in a different order
# - the canonical ordering used by the upgraded XML writer. No run
# text, paragraph, style value, or page-setup value actually changed
# anywhere in word/document.xml or word/styles.xml.
#
# So the faithful edit here is "re-save the package" - touch nothing
# in the document's content/formatting - which is exactly what
# upgrading the underlying writer library and regenerating the fixture
# amounts to from the document-model point of view.
$d = $word.ActiveDocument
$d.Save()
